$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" column (C) values for each generation-block of rows.
$ws.Range("C2:C18").Value = 8747
$ws.Range("C19:C37").Value = 8292
$ws.Range("C38:C46").Value = 8127
$ws.Range("C47:C160").Value = 7735
$ws.Range("C161:C165").Value = 7651
$ws.Range("C166:C252").Value = 7569
